$wb = $excel.ActiveWorkbook

# "Tests" sheet (variation file) - add new row for KillAllProcesses
$wsTests = $wb.Worksheets.Item("Tests")
$wsTests.Range("A10").Value = "Framework\KillAllProcesses.xaml"
$wsTests.Range("B10").Value = "Success"
$wsTests.Range("A30").Select()

# "Result" sheet - mirror the same new row
$wsResult = $wb.Worksheets.Item("Result")
$wsResult.Range("A10").Value = "Framework\KillAllProcesses.xaml"
$wsResult.Range("B10").Value = "Success"
$wsResult.Range("D18").Select()
